$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text (string) storage,
# matching the workbook author convention where Price/Volume columns
# are stored as text even when the content looks numeric.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.456.56"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "2.098.26"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "330.87"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.5224"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").Value = "0.4436"
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("D9").Value = "53.70"
$ws.Range("E9").Value = "  +17.04%  "
$ws.Range("D10").Value = "0.08947"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "1.155"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "24.48"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "2.093.09"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "6.706"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "7.709"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "96.53"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "0.00001124"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "0.06612"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "19.16"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "6.284"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "30.498.13"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "2.320"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").Value = "2.327.42"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("D27").Value = "22.30"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "2.572"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "163.59"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "132.12"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "1.195"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "1.666"
$ws.Range("E33").Value = "  +9.65%  "
$ws.Range("D34").Value = "6.164"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "3.897"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").Value = "10.45"
$ws.Range("E36").Value = "  +9.89%  "
$ws.Range("D37").Value = "0.02569"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").Value = "0.06804"
$ws.Range("D39").Value = "12.77"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").Value = "5.471"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "0.2266"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "0.6907"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "1.254"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "13.94"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "0.6357"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "2.267"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").Value = "3.629"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "1.246"
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("D50").Value = "1.247"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "81.94"
$ws.Range("E51").Value = "  -0.54%  "
